$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet activation / selection (moves tabSelected from sheet1 to sheet2,
#     bumps workbook activeTab to 1, and sets the live selection on sheet2) ---
$ws2.Activate()
[void]$ws2.Range("J12").Select()

# --- Column B width on sheet2: 17.57.. -> 24 ---
$ws2.Columns.Item(2).ColumnWidth = 23.16

# --- Scalar cell edits on sheet2 ---
$ws2.Range("J3").Value2 = 1
$ws2.Range("J10").Value2 = 5
$ws2.Range("J13").Value2 = 10
$ws2.Range("J25").Value2 = 5
$ws2.Range("J26").Value2 = 8

# --- Reorder rows 21-23: the old row23 ("gender") moves up to row21,
#     and the old rows 21-22 ("vaccinated","insurance") shift down one ---
# New row 21 ("gender", was row 23)
$ws2.Range("A21").Value2 = 22
$ws2.Range("B21").Value2 = "gender"
$ws2.Range("C21").Value2 = "int"
$ws2.Range("D21").Value2 = 1
$ws2.Range("E21").Value2 = "y"
$ws2.Range("F21").Value2 = "u"
$ws2.Range("G21").Value2 = "1=male, 2=female, 3=other"

# New row 22 ("vaccinated", was row 21)
$ws2.Range("A22").Value2 = 20
$ws2.Range("B22").Value2 = "vaccinated"
$ws2.Range("C22").Value2 = "boolean"
$ws2.Range("G22").Value2 = "0=no, 1=yes"

# New row 23 ("insurance", was row 22)
$ws2.Range("A23").Value2 = 21
$ws2.Range("B23").Value2 = "insurance"
$ws2.Range("C23").Value2 = "boolean"
$ws2.Range("G23").Value2 = "0=no, 1=yes"

# --- New row 27: duplicate of row 26 ("password_hashed") but for the
#     "password_hashed_match" field, with a value in J ---
$ws2.Range("A27").Value2 = 25
$ws2.Range("B27").Value2 = "password_hashed_match"
$ws2.Range("C27").Value2 = "string"
$ws2.Range("D27").Value2 = 20
$ws2.Range("E27").Value2 = "y"
$ws2.Range("F27").Value2 = "u"
$ws2.Range("G27").Value2 = "Password"
$ws2.Range("I27").Value2 = '"(?=.*\d)(?=.*[a-z])(?=.*[A-Z]).{8,}" '
$ws2.Range("J27").Value2 = 8

# Match the vertical-top alignment style ("s=2" in the template rows) on the
# cells of the new row that carry it.
$ws2.Range("A27").VerticalAlignment = -4160
$ws2.Range("B27").VerticalAlignment = -4160
$ws2.Range("C27").VerticalAlignment = -4160
$ws2.Range("G27").VerticalAlignment = -4160
$ws2.Range("I27").VerticalAlignment = -4160

Write-Host "done"
